$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('K2').Value = 6385
$ws.Range('K3').Value = 6591
$ws.Range('K4').Value = 1372
$ws.Range('K5').Value = 467
$ws.Range('K6').Value = 7268
$ws.Range('K7').Value = 22083

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('K6').Value = 158
$ws.Range('K7').Value = 654
$ws.Range('K8').Value = 1451
$ws.Range('K11').Value = 411
$ws.Range('K15').Value = 226
$ws.Range('K19').Value = 645
$ws.Range('K20').Value = 525
$ws.Range('K24').Value = 69
$ws.Range('K26').Value = 30
$ws.Range('K29').Value = 1191
$ws.Range('K31').Value = 246
$ws.Range('K33').Value = 964
$ws.Range('K36').Value = 281
$ws.Range('K37').Value = 750
$ws.Range('K41').Value = 155
$ws.Range('K42').Value = 816
$ws.Range('K44').Value = 183
$ws.Range('K47').Value = 150
$ws.Range('K48').Value = 276
$ws.Range('K50').Value = 104
$ws.Range('K52').Value = 584
$ws.Range('K53').Value = 283
$ws.Range('K54').Value = 434
$ws.Range('K57').Value = 82
$ws.Range('K59').Value = 40
$ws.Range('K63').Value = 63
$ws.Range('K64').Value = 139
$ws.Range('K65').Value = 516
$ws.Range('K67').Value = 866
$ws.Range('K76').Value = 302
$ws.Range('K79').Value = 556
$ws.Range('K80').Value = 78
$ws.Range('K83').Value = 474
$ws.Range('K84').Value = 178
$ws.Range('K85').Value = 1024
$ws.Range('K86').Value = 134
$ws.Range('K89').Value = 329
$ws.Range('K91').Value = 258
$ws.Range('K93').Value = 82
$ws.Range('K94').Value = 295
$ws.Range('K96').Value = 234
$ws.Range('K97').Value = 176
$ws.Range('K99').Value = 363
$ws.Range('K101').Value = 22083

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range('K3').Value = 44
$ws.Range('K6').Value = 99
$ws.Range('K7').Value = 234

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('K2').Value = 212
$ws.Range('K3').Value = 214
$ws.Range('K4').Value = 23
$ws.Range('K5').Value = 26
$ws.Range('K6').Value = 179
$ws.Range('K7').Value = 654

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range('K6').Value = 136
$ws.Range('K7').Value = 411

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range('K2').Value = 93
$ws.Range('K3').Value = 101
$ws.Range('K6').Value = 96
$ws.Range('K7').Value = 329

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('K3').Value = 354
$ws.Range('K7').Value = 1024

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('K3').Value = 167
$ws.Range('K6').Value = 212
$ws.Range('K7').Value = 584

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range('K4').Value = 13
$ws.Range('K6').Value = 119
$ws.Range('K7').Value = 283

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('K2').Value = 397
$ws.Range('K3').Value = 442
$ws.Range('K6').Value = 489
$ws.Range('K7').Value = 1451

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('K3').Value = 169
$ws.Range('K7').Value = 474

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('K3').Value = 344
$ws.Range('K6').Value = 297
$ws.Range('K7').Value = 964

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('K2').Value = 212
$ws.Range('K6').Value = 224
$ws.Range('K7').Value = 750

$ws = $wb.Worksheets.Item('New City')
$ws.Range('K6').Value = 184
$ws.Range('K7').Value = 516

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('K2').Value = 94
$ws.Range('K7').Value = 363

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range('K2').Value = 82
$ws.Range('K4').Value = 10
$ws.Range('K6').Value = 84
$ws.Range('K7').Value = 246

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('K2').Value = 237
$ws.Range('K3').Value = 317
$ws.Range('K6').Value = 244
$ws.Range('K7').Value = 866

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range('K2').Value = 60
$ws.Range('K3').Value = 70
$ws.Range('K7').Value = 178

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('K6').Value = 235
$ws.Range('K7').Value = 434

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('K2').Value = 337
$ws.Range('K3').Value = 426
$ws.Range('K7').Value = 1191

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('K3').Value = 66
$ws.Range('K7').Value = 276

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('K2').Value = 192
$ws.Range('K3').Value = 194
$ws.Range('K6').Value = 210
$ws.Range('K7').Value = 645

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range('K3').Value = 48
$ws.Range('K7').Value = 183

$ws = $wb.Worksheets.Item('River North')
$ws.Range('K2').Value = 68
$ws.Range('K3').Value = 57
$ws.Range('K7').Value = 302

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range('K2').Value = 57
$ws.Range('K7').Value = 158

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range('K2').Value = 53
$ws.Range('K7').Value = 155

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('K4').Value = 33
$ws.Range('K6').Value = 303
$ws.Range('K7').Value = 816

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range('K3').Value = 61
$ws.Range('K6').Value = 86

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range('K2').Value = 28
$ws.Range('K7').Value = 69

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range('K2').Value = 66
$ws.Range('K6').Value = 53
$ws.Range('K7').Value = 258

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('K2').Value = 185
$ws.Range('K3').Value = 180
$ws.Range('K7').Value = 556

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range('K6').Value = 53
$ws.Range('K7').Value = 139

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('K2').Value = 182
$ws.Range('K6').Value = 143
$ws.Range('K7').Value = 525

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range('K6').Value = 63
$ws.Range('K7').Value = 281

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range('K2').Value = 26
$ws.Range('K7').Value = 82

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range('K6').Value = 133
$ws.Range('K7').Value = 295

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range('K3').Value = 42
$ws.Range('K7').Value = 150

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range('K3').Value = 56
$ws.Range('K6').Value = 69
$ws.Range('K7').Value = 226

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range('K6').Value = 49
$ws.Range('K7').Value = 104

$ws = $wb.Worksheets.Item('East Village')
$ws.Range('K2').Value = 3
$ws.Range('K7').Value = 30

$ws = $wb.Worksheets.Item('Montclare')
$ws.Range('K3').Value = 12
$ws.Range('K7').Value = 40

$ws = $wb.Worksheets.Item('West Town')
$ws.Range('K2').Value = 37
$ws.Range('K7').Value = 176

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range('K4').Value = 58
$ws.Range('K7').Value = 134

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range('K5').Value = 2
$ws.Range('K7').Value = 82

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range('K6').Value = 37
$ws.Range('K7').Value = 78
